# Added errors between experimental and fit.
# Updates a few fitted-parameter values on the "Gaussian" and "Lorentzian"
# sheets, and leaves the selection/cursor where the author last clicked.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Gaussian sheet
# ---------------------------------------------------------------------
$wsGaussian = $wb.Worksheets.Item("Gaussian")
$wsGaussian.Activate()

$wsGaussian.Range("G3").Value = 0.5

$wsGaussian.Range("B4").Value = 1
$wsGaussian.Range("C4").Value = 1
$wsGaussian.Range("D4").Value = 2
$wsGaussian.Range("E4").Value = 2
$wsGaussian.Range("G4").Value = 1.8

$wsGaussian.Range("B5").Value = 100
$wsGaussian.Range("C5").Value = 0
$wsGaussian.Range("D5").Value = 600
$wsGaussian.Range("E5").Value = 1000
$wsGaussian.Range("G5").Value = 6000

$wsGaussian.Range("G6").Select()

# ---------------------------------------------------------------------
# Lorentzian sheet
# ---------------------------------------------------------------------
$wsLorentzian = $wb.Worksheets.Item("Lorentzian")
$wsLorentzian.Activate()

$wsLorentzian.Range("B5").Value = 0
$wsLorentzian.Range("C5").Value = 0
$wsLorentzian.Range("D5").Value = 100
$wsLorentzian.Range("E5").Value = 500
$wsLorentzian.Range("F5").Value = 500
$wsLorentzian.Range("G5").Value = 1000

$wsLorentzian.Range("B4").Select()
